# Apply the numeric-result updates described by the commit's OOXML diff.
# Each pair is (old run text, new run text) and is unique enough within
# the document that a whole-document Find/Replace-All is safe and will
# hit every matching occurrence (there are two near-duplicate sections
# for the scikit-learn example, hence some pairs appearing twice).

$d = $word.ActiveDocument

$replacements = @(
    @("0.974 and coefficients", "1.0075 and coefficients"),
    @("2.9594, and", "2.981, and"),
    @("2.0135", "1.9891"),
    @("0.9733 and coefficients", "1.0079 and coefficients"),
    @("2.9582, and", "2.9795, and"),
    @("2.0125", "1.9885"),
    @("array([0.97379059]) and coefficients", "array([1.05822114]) and coefficients"),
    @("array([2.97930829]), and", "array([2.98804002]), and"),
    @("array([2.00812544])", "array([2.02031678])")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()

    $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
